$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    if ($Text -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$') {
        $Cell.NumberFormat = "@"
        $Cell.Value = $Text
        $Cell.Style = "Normal"
    } else {
        $Cell.Value = $Text
    }
}

# Row-by-row updates reflecting refreshed crypto market data
Set-TextValue $ws.Range('D2') '30.335.22'
Set-TextValue $ws.Range('E2') '  +0.01%  '

Set-TextValue $ws.Range('D3') '1.935.65'
Set-TextValue $ws.Range('E3') '  +0.00%  '

Set-TextValue $ws.Range('D4') '0.9991'
Set-TextValue $ws.Range('E4') '  -0.18%  '

Set-TextValue $ws.Range('D5') '0.7631'
Set-TextValue $ws.Range('E5') '  +5.32%  '

Set-TextValue $ws.Range('D6') '248.49'
Set-TextValue $ws.Range('E6') '  -1.00%  '

Set-TextValue $ws.Range('D7') '0.9982'
Set-TextValue $ws.Range('E7') '  -0.26%  '

Set-TextValue $ws.Range('D8') '28.31'
Set-TextValue $ws.Range('E8') '  +0.90%  '

Set-TextValue $ws.Range('D9') '0.3219'
Set-TextValue $ws.Range('E9') '  -2.75%  '

Set-TextValue $ws.Range('D10') '0.07119'
Set-TextValue $ws.Range('E10') '  -1.06%  '

Set-TextValue $ws.Range('D11') '0.7909'
Set-TextValue $ws.Range('E11') '  -2.50%  '

Set-TextValue $ws.Range('D12') '0.08001'
Set-TextValue $ws.Range('E12') '  -1.13%  '

Set-TextValue $ws.Range('D13') '1.931.07'

Set-TextValue $ws.Range('D14') '5.382'
Set-TextValue $ws.Range('E14') '  -2.04%  '

Set-TextValue $ws.Range('E15') '  +0.25%  '

Set-TextValue $ws.Range('D16') '14.75'
Set-TextValue $ws.Range('E16') '  -3.36%  '

Set-TextValue $ws.Range('D17') '30.334.75'
Set-TextValue $ws.Range('E17') '  -0.05%  '

Set-TextValue $ws.Range('D18') '254.23'
Set-TextValue $ws.Range('E18') '  +1.47%  '

Set-TextValue $ws.Range('D19') '0.000008037'
Set-TextValue $ws.Range('E19') '  -3.24%  '

Set-TextValue $ws.Range('D20') '5.808'
Set-TextValue $ws.Range('E20') '  -1.85%  '

Set-TextValue $ws.Range('D21') '2.189.12'
Set-TextValue $ws.Range('E21') '  -0.05%  '

Set-TextValue $ws.Range('D22') '0.9980'
Set-TextValue $ws.Range('E22') '  -0.25%  '

Set-TextValue $ws.Range('D23') '0.9993'
Set-TextValue $ws.Range('E23') '  -0.17%  '

Set-TextValue $ws.Range('D24') '6.835'
Set-TextValue $ws.Range('E24') '  -2.29%  '

Set-TextValue $ws.Range('D25') '9.612'
Set-TextValue $ws.Range('E25') '  -1.46%  '

Set-TextValue $ws.Range('D26') '165.51'
Set-TextValue $ws.Range('E26') '  +1.11%  '

Set-TextValue $ws.Range('D27') '0.1356'
Set-TextValue $ws.Range('E27') '  +2.19%  '

Set-TextValue $ws.Range('D28') '2.322'
Set-TextValue $ws.Range('E28') '  -2.38%  '

Set-TextValue $ws.Range('D29') '19.13'
Set-TextValue $ws.Range('E29') '  -0.87%  '

Set-TextValue $ws.Range('D30') '1.374'
Set-TextValue $ws.Range('E30') '  +1.90%  '

Set-TextValue $ws.Range('E31') '  -2.59%  '

Set-TextValue $ws.Range('E32') '  -0.02%  '

Set-TextValue $ws.Range('D33') '4.152'
Set-TextValue $ws.Range('E33') '  -0.70%  '

Set-TextValue $ws.Range('D34') '0.05196'
Set-TextValue $ws.Range('E34') '  -0.12%  '

Set-TextValue $ws.Range('D35') '1.294'
Set-TextValue $ws.Range('E35') '  +0.60%  '

Set-TextValue $ws.Range('D36') '0.7548'
Set-TextValue $ws.Range('E36') '  +0.46%  '

Set-TextValue $ws.Range('D37') '2.771'
Set-TextValue $ws.Range('E37') '  +0.88%  '

Set-TextValue $ws.Range('D38') '0.01973'
Set-TextValue $ws.Range('E38') '  -0.56%  '

Set-TextValue $ws.Range('E39') '  -1.24%  '

Set-TextValue $ws.Range('D40') '78.56'
Set-TextValue $ws.Range('E40') '  -1.61%  '

Set-TextValue $ws.Range('D41') '6.455'
Set-TextValue $ws.Range('E41') '  +0.22%  '

Set-TextValue $ws.Range('D42') '0.4530'
Set-TextValue $ws.Range('E42') '  -0.24%  '

Set-TextValue $ws.Range('D43') '1.998'
Set-TextValue $ws.Range('E43') '  -1.69%  '

Set-TextValue $ws.Range('D44') '0.9988'
Set-TextValue $ws.Range('E44') '  -0.14%  '

Set-TextValue $ws.Range('D45') '0.8380'
Set-TextValue $ws.Range('E45') '  -1.29%  '

Set-TextValue $ws.Range('D46') '102.36'
Set-TextValue $ws.Range('E46') '  +0.36%  '

Set-TextValue $ws.Range('B47') 'EnergySwap'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D47') '9.832'
Set-TextValue $ws.Range('E47') '  +0.30%  '

Set-TextValue $ws.Range('B48') 'Aptos'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D48') '7.560'
Set-TextValue $ws.Range('E48') '  +0.94%  '

Set-TextValue $ws.Range('D49') '988.90'
Set-TextValue $ws.Range('E49') '  +13.01%  '

Set-TextValue $ws.Range('D50') '37.38'
Set-TextValue $ws.Range('E50') '  +1.50%  '

Set-TextValue $ws.Range('B51') 'Algorand'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range('D51') '0.1196'
Set-TextValue $ws.Range('E51') '  +4.84%  '

